# Registrations sheet rework:
#  - header renamed ("title"/"guide"/"teamMembers" -> "Project Title"/"Guide"/"USN1..4")
#  - each row's combined "USN1\nUSN2" team-members cell split into separate USN columns
#  - three new project rows appended

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Project Title', 'Guide', 'USN1', 'USN2', 'USN3', 'USN4'),
    @('Killing Dumbledore 101', 'Snape', '1DS15IS061', '1DS15IS062'),
    @('Swish and Flick', 'Flitwick', '1DS15IS063', '1DS15IS064'),
    @('Transfiguration', 'Minerva', '1DS15IS065'),
    @('Defence Against The Dark Arts', 'Snape', '1DS16IS003', '1DS16IS004'),
    @('How to bring statues to life', 'Minerva', '1DS15IS033', '1DS15IS044'),
    @("It's Levi-oh-sah, not Levio-sah", 'Hermione', '1DS16IS444'),
    @('How to betray the most evil Wizard of all time while keeping a straight face', 'Snape', '1DS16IS033')
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
